# ---------------------------------------------------------------------------
# Book1.xlsx edit: insert a new "Sheet3" (holding the former Sheet2 content)
# between Sheet1 and Sheet2, populate a new "Sheet2" with fresh notes, and
# fill in a bunch of new answers (column B, plus new columns C/D) on Sheet1.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- 1. Insert a brand-new worksheet between Sheet1 and the existing Sheet2.
#        Excel slides the existing "Sheet2" down; we rename the new blank
#        sheet to "Sheet3" so the tab order becomes Sheet1, Sheet3, Sheet2
#        (the original Sheet2 content effectively becomes "Sheet3").
$ws3 = $wb.Worksheets.Add($ws2)
$ws3.Name = "Sheet3"

# --- 2. Populate the freshly added "Sheet2" tab with new notes.
$ws2.Range("A1").Value = "cong ty dau tien lam ve sach dien tu, phat trien phan mem cho cong ty "
$ws2.Range("A2").Value = "Cong ty thu 2 lam ve offshore, cho thue nhan vien Haken"
$ws2.Range("A4").Value = "Noi dung cong viec tai 2 cong ty , da hoc hoi duoc gi o 2 cty va thich lam o cong ty nao hon "
$ws2.Range("F12").Select()

# --- 3. Fill in new interview answers on Sheet1 (column B), in the order
#        they were apparently typed in (keeps shared-string / style table
#        ordering consistent with the source file).
$ws1.Range("B6").Value = "Tiep xuc voi khach hang, thiet ke, viet tai lieu, quan ly du an, hoc ky thuat moi"

$ws1.Range("B10").Value = "仕事を進めるにあたって、コミュニケーションを大事にする人や効率を重視する人など、大事にするものは人それぞれです。"

$ws1.Range("B12").Value = "Duoc khen,cong viec dien ra suon se thi se co cam giac dat duoc thanh tuu. Tien, thuong"

$ws1.Range("B13").Value = "Toi chua gap nguoi khong ua bao gio. Neu ma gap thi co le toi se co gang giao tiep bang van ban, hoac bay to y kien 1 cach can than va chi tiet hon de tranh conflic khi lam viec"

# B14 is a hyperlink cell -> give it the link + wrap text (this establishes
# the new "hyperlink + wrap" cell style used by the other link cells below).
$ws1.Range("B14").Value = "https://www.theport.jp/portcareer/article/31918/"
$ws1.Hyperlinks.Add($ws1.Range("B14"), "https://www.theport.jp/portcareer/article/31918/") | Out-Null
$ws1.Range("B14").WrapText = $true

$ws1.Range("B15").Value = "het minh vi cong viec, khong ngai lam them gio, giao tiep, mong muon hoc cong nghe moi, dong luc -> Ke ra vi du thuc te, lam them gio khi lam du an o VN, nhiet huyet, cam giac thanh cong khi minh hoan tot mot cong viec j do"

$ws1.Range("B16").Value = "talkative, lam luc noi hoi nhieu"

$ws1.Range("B19").Value = "co hoi lam viec tai nhat, co hoc hoi phong cach lam viec,tieng nhat va sach dien tu cung thu vi"

# B21 / B22: more hyperlinked answers.
$ws1.Range("B21").Value = "https://media.bizreach.biz/member_survey1/#:~:text=%E8%BB%A2%E8%81%B7%E8%BB%B8%E3%81%AF%E5%A4%A7%E3%81%8D%E3%81%8F%E3%80%81%E3%80%8C%E4%BA%BA,%E4%BB%98%E3%81%91%E3%81%AB%E3%82%82%E3%81%A4%E3%81%AA%E3%81%8C%E3%82%8A%E3%81%BE%E3%81%99%E3%80%82"
$ws1.Hyperlinks.Add($ws1.Range("B21"), "https://media.bizreach.biz/member_survey1/", ":~:text=%E8%BB%A2%E8%81%B7%E8%BB%B8%E3%81%AF%E5%A4%A7%E3%81%8D%E3%81%8F%E3%80%81%E3%80%8C%E4%BA%BA,%E4%BB%98%E3%81%91%E3%81%AB%E3%82%82%E3%81%A4%E3%81%AA%E3%81%8C%E3%82%8A%E3%81%BE%E3%81%99%E3%80%82") | Out-Null
$ws1.Range("B21").WrapText = $true

$ws1.Range("B22").Value = "https://ourly.jp/joined-company_gap/"
$ws1.Hyperlinks.Add($ws1.Range("B22"), "https://ourly.jp/joined-company_gap/") | Out-Null
$ws1.Range("B22").WrapText = $true

# Highlight the "gap after joining" question in red now that the link style
# above already exists (keeps the new style table order matching source).
$ws1.Range("A22").Font.Color = 255

# D1: new hyperlinked column header-ish note.
$ws1.Range("D1").Value = "https://one-group.jp/tenshoku/tenshoku-kikkake/"
$ws1.Hyperlinks.Add($ws1.Range("D1"), "https://one-group.jp/tenshoku/tenshoku-kikkake/") | Out-Null

# C1: new column with a longer Japanese note (plain wrap style, no link).
$ws1.Range("C1").Value = "スキルアップがしたかったから`n　　　転職先ではスキルが評価され、プロジェクトリーダーなど責任ある仕事を任されることもあります。`n会社の将来に不安を覚えたから`n          給料や待遇の水準が良いというメリットもあります。"

# B5 gets replaced with a fuller answer (used to just be the stray "1, ").
$ws1.Range("B5").Value = "1, Muon lam Brse`n【PREP法の構成】`nP=point（結論）`nR=reason（理由）`nE=example（具体例）`nP=point（結論）"

# B7: hyperlinked multi-line answer (two URLs stacked, link points at the
# first one and shows it as the display text).
$ws1.Range("B7").Value = "https://career-ch.com/column/3080`nhttps://career-ch.com/column/27916"
$ws1.Hyperlinks.Add($ws1.Range("B7"), "https://career-ch.com/column/3080", "", "", "https://career-ch.com/column/3080") | Out-Null
$ws1.Range("B7").WrapText = $true

# --- 4. Sheet1 layout tweaks: new column widths for the added C/D columns,
#        zoom + scroll position, active selection.
$ws1.Columns.Item(3).ColumnWidth = 62.58
$ws1.Columns.Item(4).ColumnWidth = 30.17

$ws1.Activate()
$ws1.Range("A16").Select()
$excel.ActiveWindow.Zoom = 88
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "edit complete"
